$wb = $excel.ActiveWorkbook

# Overview sheet: G2 "Latest HO Xliff Generate Date"
$ws = $wb.Worksheets.Item("Overview")
$ws.Range("G2").Value = "2016-09-05 07:15:18"

# zh-cn sheet: H2 "Correspond Handoff Datetime", K2 "Correspond Handback DateTime"
$ws = $wb.Worksheets.Item("zh-cn")
$ws.Range("H2").Value = "2016-09-05 07:15:06"
$ws.Range("K2").Value = "2016-09-05 07:15:33"

# de-de sheet: H2 "Correspond Handoff Datetime", K2 "Correspond Handback DateTime"
$ws = $wb.Worksheets.Item("de-de")
$ws.Range("H2").Value = "2016-09-05 07:15:18"
$ws.Range("K2").Value = "2016-09-05 07:15:41"
